$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.9633038705127712
$ws0.Range("C2").Value = -0.2889435717831422
$ws0.Range("B3").Value = 0.8870168358053145
$ws0.Range("C3").Value = 0.3890749122423709
$ws0.Range("B4").Value = -0.2053854463661027
$ws0.Range("C4").Value = -0.3121330545868284

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.3778835960970965
$ws1.Range("C2").Value = -0.6208671227570343
$ws1.Range("B3").Value = -0.7967644956872764
$ws1.Range("C3").Value = 0.9125765515905209
$ws1.Range("B4").Value = -0.2199469702406884
$ws1.Range("C4").Value = -0.3915850869071731
